$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1/IF1 with the same style as H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-33: I = 1, J = same value as H
For ($r = 2; $r -le 33; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}

# Row 34 is the exception: I = 2, J = 2
$ws.Cells.Item(34, 9).Value = 2
$ws.Cells.Item(34, 10).Value = 2
